# EvenStarFarm.xlsx edit script
# Adds the Solver constraints block (rows 39-51, column G, updated decision
# variables, etc.) as described by the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Decision-variable cells B26:D33 — fill in the optimal solution values
#    that used to be blank.
# ---------------------------------------------------------------------------
$ws.Range("B26").Value = 416
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 608
$ws.Range("D27").Value = 0

$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 167

$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 86
$ws.Range("D29").Value = 0

$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 72
$ws.Range("D30").Value = 0

$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 251

$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 107

$ws.Range("B33").Value = 58
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 75

# ---------------------------------------------------------------------------
# 2) B9 becomes a small formula (10 + the original 76) instead of a literal.
# ---------------------------------------------------------------------------
$ws.Range("B9").Formula = "=10+76"

# ---------------------------------------------------------------------------
# 3) New helper column G6:G13 (per-produce profit contribution, literal
#    cached numbers - no formula in the source file either).
# ---------------------------------------------------------------------------
$ws.Range("G6").Value = 713.47119327733526
$ws.Range("G7").Value = 331.49
$ws.Range("G8").Value = 184
$ws.Range("G9").Value = 331.49
$ws.Range("G10").Value = 552.48
$ws.Range("G11").Value = 565.47
$ws.Range("G12").Value = 673.47
$ws.Range("G13").Value = 673.47

# ---------------------------------------------------------------------------
# 4) C36 gains a small check formula against A37.
# ---------------------------------------------------------------------------
$ws.Range("C36").Formula = "=A37-49956.3917680672"

# ---------------------------------------------------------------------------
# 5) New "constraints" block, rows 39-51.
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "constraints"
$ws.Range("A39").Font.Bold = $true
$ws.Range("A39").Font.Size = 12

$ws.Range("A41:A48").Merge()
$ws.Range("A41").Value = "available"

$ws.Range("C40").Value = "LHS"
$ws.Range("D40").Value = "sign"
$ws.Range("E40").Value = "RHS"
$ws.Range("C40:E40").Font.Bold = $false
$ws.Range("C40:E40").HorizontalAlignment = -4108  # xlCenter

$produceDataRows = 6..13
$decisionRows = 26..33
$constraintRows = 41..48
for ($i = 0; $i -lt $decisionRows.Count; $i++) {
    $dr = $decisionRows[$i]
    $pd = $produceDataRows[$i]
    $cr = $constraintRows[$i]
    $ws.Cells.Item($cr, 2).Value = $ws.Cells.Item($dr, 1).Value2
    $ws.Cells.Item($cr, 3).Formula = "=SUM(B" + $dr + ":D" + $dr + ")"
    $ws.Cells.Item($cr, 4).Value = "<="
    $ws.Cells.Item($cr, 5).Formula = "=B" + $pd
}

$ws.Range("A49:B49").Merge()
$ws.Range("A49").Value = "truck capacity"
$ws.Range("C49").Formula = "=SUM(D26:D33)"
$ws.Range("D49").Value = "<="
$ws.Range("E49").Formula = "=600"

$ws.Range("A50:B50").Merge()
$ws.Range("A50").Value = "restaurants limit"
$ws.Range("C50").Formula = "=SUM(B26:B33)/119"
$ws.Range("D50").Value = "<="
$ws.Range("E50").Value = 20

$ws.Range("A51:B51").Merge()
$ws.Range("A51").Value = "CSA customers limit"
$ws.Range("C51").Formula = "=SUMPRODUCT(C26:C33,D6:D13)/400"
$ws.Range("D51").Value = "<="
$ws.Range("E51").Formula = "=90"

# Light styling for the new block: thin borders around the numeric table,
# bold label font, and a number format on the LHS/RHS columns. Kept to the
# cells that already carry content so we don't mint empty styled cells that
# don't exist in the target sheet (e.g. B40).
$table = $ws.Range("C40:E51")
$table.Borders.Item(7).LineStyle = 1
$table.Borders.Item(8).LineStyle = 1
$table.Borders.Item(9).LineStyle = 1
$table.Borders.Item(10).LineStyle = 1
$table.Borders.Item(11).LineStyle = 1
$table.Borders.Item(12).LineStyle = 1

$bcol = $ws.Range("B41:B51")
$bcol.Borders.Item(7).LineStyle = 1
$bcol.Borders.Item(8).LineStyle = 1
$bcol.Borders.Item(9).LineStyle = 1
$bcol.Borders.Item(10).LineStyle = 1
$bcol.Borders.Item(11).LineStyle = 1
$bcol.Borders.Item(12).LineStyle = 1

$ws.Range("C41:C48").NumberFormat = "0"
$ws.Range("E41:E48").NumberFormat = "0"
$ws.Range("C49:C51").NumberFormat = "#,##0.00_);[Red]\(#,##0.00\)"
$ws.Range("E49:E51").NumberFormat = "0"

foreach ($addr in @("A41", "A49", "A50", "A51")) {
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 6) Solver parameters (hidden workbook-scoped / sheet-scoped defined names)
#    recorded by Excel's Solver add-in when the model above was solved.
# ---------------------------------------------------------------------------
$solverNames = @(
    @("solver_adj",  "=Sheet1!`$B`$26:`$D`$33"),
    @("solver_cvg",  "0.0001"),
    @("solver_drv",  "1"),
    @("solver_eng",  "2"),
    @("solver_est",  "1"),
    @("solver_itr",  "2147483647"),
    @("solver_lhs1", "=Sheet1!`$C`$41:`$C`$51"),
    @("solver_mip",  "2147483647"),
    @("solver_mni",  "30"),
    @("solver_mrt",  "0.075"),
    @("solver_msl",  "2"),
    @("solver_neg",  "1"),
    @("solver_nod",  "2147483647"),
    @("solver_num",  "1"),
    @("solver_nwt",  "1"),
    @("solver_opt",  "=Sheet1!`$A`$37"),
    @("solver_pre",  "0.000001"),
    @("solver_rbv",  "1"),
    @("solver_rel1", "1"),
    @("solver_rhs1", "=Sheet1!`$E`$41:`$E`$51"),
    @("solver_rlx",  "2"),
    @("solver_rsd",  "0"),
    @("solver_scl",  "1"),
    @("solver_sho",  "2"),
    @("solver_ssz",  "100"),
    @("solver_tim",  "2147483647"),
    @("solver_tol",  "0.01"),
    @("solver_typ",  "1"),
    @("solver_val",  "0"),
    @("solver_ver",  "3")
)

foreach ($pair in $solverNames) {
    $n = $ws.Names.Add($pair[0], $pair[1])
    $n.Visible = $false
}

Write-Output "EvenStarFarm solver block applied"
